$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Restructure sheet1: insert 3 blank rows at row 5 (pushes old row5.. down
#    by 3), then insert 1 more blank row at row 9 (pushes old row6(now at
#    row9)... down by 1 more). This reproduces the row layout of the target
#    sheet while preserving formatting/formulas of all the untouched rows.
# ---------------------------------------------------------------------------
$ws1.Rows("5:7").Insert()
$ws1.Rows("9:9").Insert()

# ---------------------------------------------------------------------------
# 2. Row 4 - fill in the new "Binding posts" part (existing blank template
#    row already carries the correct styles for B4/C4/D4/E4/G4).
# ---------------------------------------------------------------------------
$ws1.Range("A4").Value = "Binding posts"
$ws1.Range("A2").Copy()
$ws1.Range("A4").PasteSpecial(-4122)
$ws1.Range("A4").Value = "Binding posts"

$ws1.Range("B4").Value = "5-way dual posts"
$ws1.Range("C4").Value = "farnell"
$ws1.Range("D4").Value = 2.56
$ws1.Range("E4").Formula = "=D4*1.21"

$ws1.Range("F4").Value = 1
$ws1.Range("F2").Copy()
$ws1.Range("F4").PasteSpecial(-4122)
$ws1.Range("F4").Value = 1

$ws1.Range("G4").Formula = "=F4*E4"

Write-Host "row4 done"
